$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 229, pushing the existing rows 229-246 down to 230-247.
$ws.Range("A229").EntireRow.Insert()

# Populate the newly inserted row 229 with the new weekly record.
$ws.Range("A229").Value = 5
$ws.Range("B229").Value = "Macroferia Regional de Talca"
$ws.Range("C229").Value = "Maule"
$ws.Range("D229").Value = 44714
$ws.Range("E229").Value = 7
$ws.Range("F229").Value = 100112008
$ws.Range("G229").Value = "Coliflor"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 2000
$ws.Range("K229").Value = 1100
$ws.Range("L229").Value = 1100
$ws.Range("M229").Value = 1100
$ws.Range("N229").Value = "$/unidad"
$ws.Range("O229").Value = "Región del Maule"
$ws.Range("P229").Value = 1100
$ws.Range("Q229").Value = 1
$ws.Range("R229").Value = "Hortaliza"
